# Add data for 2022-10-05
# - Rename the sheet/header to reflect the new "through" date (09-25 -> 09-27)
# - Update September 2022 figure (I10) and the Total 2022 figure (I14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Through 2022-09-27"

# Update the column header string (I1) which reads "2022 (through 09-25)"
$ws.Range("I1").Value = "2022 (through 09-27)"

# Update September 2022 value
$ws.Range("I10").Value = 127

# Update Total 2022 value
$ws.Range("I14").Value = 1262
